$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 6
$ws.Range("F3").Value = 5
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = -5
$ws.Range("F7").Value = -13
